# Update column G (K - strikeouts) values per regenerated save_data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 1
    4 = 1
    5 = 0
    6 = 0
    7 = 1
    9 = 2
    10 = 1
    11 = 3
    12 = 0
    13 = 0
    14 = 1
    16 = 3
    17 = 1
    18 = 3
    19 = 1
    20 = 0
    21 = 2
    22 = 1
    23 = 3
    24 = 3
    25 = 3
    26 = 3
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 3
    32 = 1
    33 = 0
    34 = 2
    35 = 3
    36 = 0
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 1
    42 = 0
    43 = 2
    44 = 1
    45 = 2
    46 = 2
    47 = 1
    48 = 1
    49 = 2
    50 = 2
    51 = 0
    52 = 0
    53 = 1
    54 = 0
    55 = 0
    56 = 1
    57 = 1
    58 = 2
    61 = 2
    62 = 2
    64 = 0
    65 = 2
    66 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

